# Natmi following Dr Hou advice
# Update the LR-pair stats for rows 2-10 (columns E..T) to reflect the
# recomputed ligand/receptor-expressing cell counts (1 -> 3) and the
# resulting downstream expression/specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3, 1, 0.747119, 2.241357, 0.03096954854571248, 0.03096954854571248, 3, 1, 2.324764666666667, 6.974294, 0.04473923998638302, 0.04473923998638301, 1.736875852995333, 15.631882676958, 0.00138555406465657, 0.00138555406465657)
    3  = @(3, 1, 0.747119, 2.241357, 0.03096954854571248, 0.03096954854571248, 3, 1, 19.27491966666667, 57.824759, 0.3709387315842666, 0.3709387315842665, 14.40065870644033, 129.605928357963, 0.01148780505528395, 0.01148780505528395)
    4  = @(3, 1, 0.747119, 2.241357, 0.03096954854571248, 0.03096954854571248, 3, 1, 30.36285833333334, 91.08857500000001, 0.5843220284293504, 0.5843220284293504, 22.68466835514167, 204.162015196275, 0.01809618942577195, 0.01809618942577195)
    5  = @(3, 1, 19.74619233333334, 59.23857700000001, 0.8185184181638298, 0.8185184181638298, 3, 1, 2.324764666666667, 6.974294, 0.04473923998638302, 0.04473923998638301, 45.90525023773757, 413.1472521396381, 0.03661989194350618, 0.03661989194350618)
    6  = @(3, 1, 19.74619233333334, 59.23857700000001, 0.8185184181638298, 0.8185184181638298, 3, 1, 19.27491966666667, 57.824759, 0.3709387315842666, 0.3709387315842665, 380.6062709475493, 3425.456438527943, 0.3036201838120514, 0.3036201838120513)
    7  = @(3, 1, 19.74619233333334, 59.23857700000001, 0.8185184181638298, 0.8185184181638298, 3, 1, 30.36285833333334, 91.08857500000001, 0.5843220284293504, 0.5843220284293504, 599.550840439753, 5395.957563957776, 0.4782783424082723, 0.4782783424082723)
    8  = @(3, 1, 3.630999, 10.892997, 0.1505120332904577, 0.1505120332904577, 3, 1, 2.324764666666667, 6.974294, 0.04473923998638302, 0.04473923998638301, 8.441218179902, 75.970963619118, 0.006733793978220256, 0.006733793978220256)
    9  = @(3, 1, 3.630999, 10.892997, 0.1505120332904577, 0.1505120332904577, 3, 1, 19.27491966666667, 57.824759, 0.3709387315842666, 0.3709387315842665, 69.98721403474698, 629.884926312723, 0.05583074271693128, 0.05583074271693128)
    10 = @(3, 1, 3.630999, 10.892997, 0.1505120332904577, 0.1505120332904577, 3, 1, 30.36285833333334, 91.08857500000001, 0.5843220284293504, 0.5843220284293504, 110.247508245475, 992.227574209275, 0.08794749659530615, 0.08794749659530617)
}

# Columns E (5) through T (20)
$firstCol = 5

foreach ($rowNum in $data.Keys | Sort-Object) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $firstCol + $i
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}
